$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.169.95"
$ws.Range("E2").Value = "  +4.36%  "
$ws.Range("D3").Value = "1.907.11"
$ws.Range("E3").Value = "  +5.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5093"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.16"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3028"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06818"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "1.905.77"
$ws.Range("E11").Value = "  +5.03%  "
$ws.Range("E12").Value = "  +3.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07334"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6940"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.923"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.31%  "
$ws.Range("D17").Value = "30.154.85"
$ws.Range("E17").Value = "  +4.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008266"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +13.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9986"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.91%  "
$ws.Range("D21").Value = "2.151.08"
$ws.Range("E21").Value = "  +5.30%  "
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.834"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.762"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.394"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "148.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "134.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.011"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.400"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.311"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08906"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.008"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05178"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.152"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7234"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.686"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.45%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.821"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.74%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.301"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9635"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01691"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.106"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4334"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "105.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.73%  "
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.700"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1284"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05765"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.417"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3831"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.01%  "
